$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p112r_4</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p112r_4</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p112v_1</id>", $false, $false, $false, $false, $false, $true, 1, $false, "<id>p112v_1</id>", 2) | Out-Null
